$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.469.82'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').Value = '2.255.77'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.18'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.437'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0952'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -7.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.93'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.30'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').Value = '2.590.51'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.92'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.04'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.821'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '2.265.82'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').Value = '43.384.34'
$ws.Range('E19').Value = '  -1.33%  '
$ws.Range('D20').Value = '0.0₃0965'
$ws.Range('E20').Value = '  -4.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.07'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '246.70'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +19.39%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('E27').Value = '  -2.36%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '173.57'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.59'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +4.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.42'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.13%  '
$ws.Range('E32').Value = '  -4.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.124'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('E34').Value = '  +4.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0676'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.89'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('E37').Value = '  -5.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.40'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.59%  '
$ws.Range('E39').Value = '  -1.56%  '
$ws.Range('E40').Value = '  -3.24%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.77'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.50'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.07'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.41'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0938'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.54%  '
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'Celestia'
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.09'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000205'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').Value = '1.423.44'
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.25'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.95%  '
